# Fix a bug in the slicing assumption set (portfolio_small.xlsx):
# the DATE_OF_DISABLEMENT column (L) carried a stray hard-coded disablement
# date on the one sample policy, even though that record is not actually
# disabled. Clear it back to the sheet's normal "unset" placeholder and give
# it the same date display format used elsewhere in the sheet (instead of the
# one-off built-in date format that was only ever used by this cell), and
# tidy up the selection that had been accidentally left spanning the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("L2")
$cell.NumberFormat = "dd/mm/yyyy"
$cell.Value = 1

# Selection had drifted to the whole row (A3:L3); restore it to just L3.
$ws.Range("L3").Select()
